$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{Row=2; D=10; E=9.5; F=11.40999984741211; G=9.109999656677246; H=2274261581}
    @{Row=3; D=10; E=9.5; F=11.40999984741211; G=9.109999656677246; H=2274261581}
    @{Row=4; D=10; E=9.5; F=11.40999984741211; G=9.109999656677246; H=2274261581}
    @{Row=5; D=10; E=9.5; F=11.40999984741211; G=9.109999656677246; H=2274261581}
    @{Row=6; D=10; E=9.5; F=11.40999984741211; G=9.109999656677246; H=2274261581}
    @{Row=7; D=9.689999580383301; E=10.13000011444092; F=11.10000038146973; G=8.899999618530273; H=2274261581}
    @{Row=8; D=23.90999984741211; E=35.18000030517578; F=45; G=22.5; H=2274261581}
    @{Row=9; D=23.95000076293945; E=23.04000091552734; F=26.20000076293945; G=21.23999977111816; H=2274261581}
    @{Row=10; D=26.10000038146973; E=21.70999908447266; F=26.1299991607666; G=20.54999923706055; H=2274261581}
    @{Row=11; D=24.20999908447266; E=25.8799991607666; F=26.75; G=23.02000045776367; H=2274261581}
    @{Row=12; D=18.36000061035156; E=13.71000003814697; F=18.84000015258789; G=11.75; H=2274261581}
    @{Row=13; D=13.77000045776367; E=10.39999961853027; F=14.85999965667725; G=10.35999965667725; H=2274261581}
    @{Row=14; D=9.109999656677246; E=10.35000038146973; F=10.76000022888184; G=8.649999618530273; H=2274261581}
    @{Row=15; D=8.145000457763672; E=8.789999961853027; F=9.050000190734863; G=7.53000020980835; H=2274261581}
    @{Row=16; D=6.579999923706055; E=7.78000020980835; F=7.820000171661377; G=6.090000152587891; H=2274261581}
    @{Row=17; D=8.369999885559082; E=7.75; F=9.039999961853027; G=7.619999885559082; H=2274261581}
    @{Row=18; D=15.48999977111816; E=19.84000015258789; F=19.9950008392334; G=14.61999988555908; H=2274261581}
    @{Row=19; D=16.03000068664551; E=14.80000019073486; F=18.44000053405762; G=14.55000019073486; H=2274261581}
    @{Row=20; D=16.95000076293945; E=16.09000015258789; F=18.35000038146973; G=15.66399955749512; H=2274261581}
    @{Row=21; D=22.97999954223633; E=21.96999931335449; F=24.10000038146973; G=20.32999992370605; H=2274261581}
    @{Row=22; D=25.47999954223633; E=26.88999938964844; F=29.82999992370605; G=25.14200019836425; H=2274261581}
    @{Row=23; D=37.20999908447266; E=41.56000137329102; F=45.13999938964844; G=36.04999923706055; H=2274261581}
    @{Row=24; D=76.19999694824219; E=82.48999786376953; F=85.22000122070312; G=63.40000152587891; H=2274261581}
    @{Row=25; D=83.88999938964844; E=118.4400024414062; F=118.7799987792969; G=66.12000274658203; H=2274261581}
    @{Row=26; D=135.2700042724609; E=158.3500061035156; F=160.8899993896484; G=128.5099945068359; H=2274261581}
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value = $row.D   # D: open_price
    $ws.Cells.Item($r, 5).Value = $row.E   # E: close_price
    $ws.Cells.Item($r, 6).Value = $row.F   # F: high_price
    $ws.Cells.Item($r, 7).Value = $row.G   # G: low_price
    $ws.Cells.Item($r, 8).Value = $row.H   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "PLTR"   # I: fixed_ticker
}
